$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 159, pushing existing rows 159..261 down to 160..262
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with the new data point
$ws.Cells.Item(159, 1).Value = 8
$ws.Cells.Item(159, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(159, 3).Value = "Coquimbo"
$ws.Cells.Item(159, 4).Value = "2022-04-29"
$ws.Cells.Item(159, 5).Value = 4
$ws.Cells.Item(159, 6).Value = 100112003
$ws.Cells.Item(159, 7).Value = "Ajo"
$ws.Cells.Item(159, 8).Value = "Chino"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 560
$ws.Cells.Item(159, 11).Value = 19000
$ws.Cells.Item(159, 12).Value = 20000
$ws.Cells.Item(159, 13).Value = 19500
$ws.Cells.Item(159, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(159, 15).Value = "China"
$ws.Cells.Item(159, 16).Value = 1950
$ws.Cells.Item(159, 17).Value = 10
$ws.Cells.Item(159, 18).Value = "Hortaliza"
